$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column D (Sample_ID), shifting
# Sample_ID and relative_file_path to the right.
$ws.Columns("D").Insert()

# New header + value for the inserted "is_normal_for_donor" column.
$ws.Range("D1").Value = "is_normal_for_donor"
$ws.Range("D2").Value = "Y"

# The new column keeps a manual (non-autofit) width matching its
# left neighbour rather than inheriting the old column D's bestFit width.
$ws.Columns("D").ColumnWidth = 8.25

# Update the selection to match the post-edit cursor position.
$ws.Range("D4").Select()
